# feat: add 2022-Q1 data
#
# The sheet that used to be "总计" (3rd sheet) is reused/renamed to become
# "2022-Q1" (keeping its identity / sheetId), and a brand-new worksheet is
# appended right after it to become the new "总计" (gets a fresh sheetId).
# This mirrors exactly how the target workbook.xml numbers sheetId/r:id for
# the two sheets.

$wb = $excel.ActiveWorkbook

# Sheet with the same column layout ("基金代码"/"基金名称"/... holdings
# table) used purely as a formatting template, so copied cells reuse the
# existing style entries (bold+bordered header / index column) instead of
# creating new ones.
$template = $wb.Worksheets.Item(2)

$quarterSheet = $wb.Worksheets.Item(3)
$quarterSheet.Cells.Clear()
$quarterSheet.Name = "2022-Q1"

# ---- Header row ----
$template.Range("A1:H1").Copy($quarterSheet.Range("A1:H1"))
$quarterSheet.Range("A1").Clear()
$quarterSheet.Range("B1").Value = "基金代码"
$quarterSheet.Range("C1").Value = "基金名称"
$quarterSheet.Range("D1").Value = "基金规模"
$quarterSheet.Range("E1").Value = "股票总仓位"
$quarterSheet.Range("F1").Value = "仓位占比"
$quarterSheet.Range("G1").Value = "持有市值(亿元)"
$quarterSheet.Range("H1").Value = "仓位排名"

# ---- Data rows ----
# code, name, scale, stockPosition, positionRatio, marketValue, positionRank
$quarterRows = @(
    @("001917", "招商量化精选股票A", "2.33", "94.20", "1.45", "0.0338", 4),
    @("004194", "招商中证1000指数增强A", "1.76", "94.40", "1.12", "0.0197", 2),
    @("007950", "招商量化精选股票C", "0.56", "94.20", "1.45", "0.0081", 4),
    @("004195", "招商中证1000指数增强C", "0.68", "94.40", "1.12", "0.0076", 2)
)

for ($i = 0; $i -lt $quarterRows.Count; $i++) {
    $r = $i + 2
    $template.Range("A2:H2").Copy($quarterSheet.Range("A$r`:H$r"))

    $data = $quarterRows[$i]

    $quarterSheet.Range("A$r").Value = $i

    # Force text storage (no leading-zero / decimal-string mangling) for the
    # numeric-looking values without otherwise altering their format.
    $quarterSheet.Range("B$r").NumberFormat = "@"
    $quarterSheet.Range("B$r").Value = $data[0]

    $quarterSheet.Range("C$r").Value = $data[1]

    $quarterSheet.Range("D$r").NumberFormat = "@"
    $quarterSheet.Range("D$r").Value = $data[2]

    $quarterSheet.Range("E$r").NumberFormat = "@"
    $quarterSheet.Range("E$r").Value = $data[3]

    $quarterSheet.Range("F$r").NumberFormat = "@"
    $quarterSheet.Range("F$r").Value = $data[4]

    $quarterSheet.Range("G$r").NumberFormat = "@"
    $quarterSheet.Range("G$r").Value = $data[5]

    $quarterSheet.Range("H$r").Value = $data[6]
}

# ---- New "总计" sheet, right after "2022-Q1" ----
$totalSheet = $wb.Worksheets.Add($null, $quarterSheet)
$totalSheet.Name = "总计"

$template.Range("A1:D1").Copy($totalSheet.Range("A1:D1"))
$totalSheet.Range("A1").Clear()
$totalSheet.Range("B1").Value = "日期"
$totalSheet.Range("C1").Value = "持有数量(只)"
$totalSheet.Range("D1").Value = "持有市值(亿元)"

# date, count, marketValue
$totalRows = @(
    @("2022-Q1", 4, 0.07),
    @("2021-Q4", 4, 0.05),
    @("2020-Q4", 1, 0.6)
)

for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $r = $i + 2
    $template.Range("A2:D2").Copy($totalSheet.Range("A$r`:D$r"))

    $data = $totalRows[$i]

    $totalSheet.Range("A$r").Value = $i
    $totalSheet.Range("B$r").Value = $data[0]
    $totalSheet.Range("C$r").Value = $data[1]
    $totalSheet.Range("D$r").Value = $data[2]
}
